$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.777.26"
$ws.Cells.Item(2, 5).Value = "  -6.67%  "

$ws.Cells.Item(3, 4).Value = "3.527.06"
$ws.Cells.Item(3, 5).Value = "  -2.91%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.44%  "

$ws.Cells.Item(5, 4).Value = "'395.46"
$ws.Cells.Item(5, 5).Value = "  -5.92%  "

$ws.Cells.Item(6, 4).Value = "'124.14"
$ws.Cells.Item(6, 5).Value = "  -5.34%  "

$ws.Cells.Item(7, 4).Value = "3.518.71"
$ws.Cells.Item(7, 5).Value = "  -2.95%  "

$ws.Cells.Item(8, 4).Value = "'0.592"
$ws.Cells.Item(8, 5).Value = "  -8.99%  "

$ws.Cells.Item(9, 4).Value = "'0.999"
$ws.Cells.Item(9, 5).Value = "  -0.06%  "

$ws.Cells.Item(10, 4).Value = "'0.684"
$ws.Cells.Item(10, 5).Value = "  -11.96%  "

$ws.Cells.Item(11, 4).Value = "'0.154"
$ws.Cells.Item(11, 5).Value = "  -15.97%  "

$ws.Cells.Item(12, 4).Value = "'0.0000346"
$ws.Cells.Item(12, 5).Value = "  -1.42%  "

$ws.Cells.Item(13, 4).Value = "'39.23"
$ws.Cells.Item(13, 5).Value = "  -8.53%  "

$ws.Cells.Item(14, 4).Value = "4.069.97"
$ws.Cells.Item(14, 5).Value = "  -3.06%  "

$ws.Cells.Item(15, 4).Value = "'9.27"
$ws.Cells.Item(15, 5).Value = "  -7.46%  "

$ws.Cells.Item(16, 5).Value = "  -3.15%  "

$ws.Cells.Item(17, 4).Value = "3.511.15"
$ws.Cells.Item(17, 5).Value = "  -1.94%  "

$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).Value = "'18.85"
$ws.Cells.Item(18, 5).Value = "  -8.03%  "

$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(19, 4).Value = "'12.72"
$ws.Cells.Item(19, 5).Value = "  +1.79%  "

$ws.Cells.Item(20, 4).Value = "63.677.68"
$ws.Cells.Item(20, 5).Value = "  -6.66%  "

$ws.Cells.Item(21, 4).Value = "'1.03"
$ws.Cells.Item(21, 5).Value = "  -10.17%  "

$ws.Cells.Item(22, 4).Value = "'397.85"
$ws.Cells.Item(22, 5).Value = "  -14.61%  "

$ws.Cells.Item(23, 4).Value = "'13.95"
$ws.Cells.Item(23, 5).Value = "  +4.37%  "

$ws.Cells.Item(24, 4).Value = "'81.66"
$ws.Cells.Item(24, 5).Value = "  -8.49%  "

$ws.Cells.Item(25, 5).Value = "  -8.07%  "

$ws.Cells.Item(26, 4).Value = "'33.92"
$ws.Cells.Item(26, 5).Value = "  -5.84%  "

$ws.Cells.Item(27, 4).Value = "'5.24"
$ws.Cells.Item(27, 5).Value = "  +7.30%  "

$ws.Cells.Item(28, 4).Value = "'3.00"
$ws.Cells.Item(28, 5).Value = "  -11.07%  "

$ws.Cells.Item(29, 4).Value = "'8.86"
$ws.Cells.Item(29, 5).Value = "  -12.80%  "

$ws.Cells.Item(30, 4).Value = "'11.97"
$ws.Cells.Item(30, 5).Value = "  -3.51%  "

$ws.Cells.Item(31, 4).Value = "'2.56"
$ws.Cells.Item(31, 5).Value = "  -7.98%  "

$ws.Cells.Item(32, 4).Value = "'0.111"
$ws.Cells.Item(32, 5).Value = "  -5.60%  "

$ws.Cells.Item(33, 4).Value = "'6.87"
$ws.Cells.Item(33, 5).Value = "  -7.18%  "

$ws.Cells.Item(34, 4).Value = "'0.150"
$ws.Cells.Item(34, 5).Value = "  -7.23%  "

$ws.Cells.Item(35, 5).Value = "  +0.13%  "

$ws.Cells.Item(36, 4).Value = "'36.86"
$ws.Cells.Item(36, 5).Value = "  -9.53%  "

$ws.Cells.Item(37, 4).Value = "'54.05"
$ws.Cells.Item(37, 5).Value = "  -4.92%  "

$ws.Cells.Item(38, 4).Value = "'0.0440"
$ws.Cells.Item(38, 5).Value = "  -11.29%  "

$ws.Cells.Item(39, 4).Value = "'0.997"
$ws.Cells.Item(39, 5).Value = "  -0.25%  "

$ws.Cells.Item(40, 4).Value = "'2.80"
$ws.Cells.Item(40, 5).Value = "  +19.79%  "

$ws.Cells.Item(41, 4).Value = "0.0₃0639"
$ws.Cells.Item(41, 5).Value = "  -10.47%  "

$ws.Cells.Item(42, 4).Value = "'0.132"
$ws.Cells.Item(42, 5).Value = "  -9.72%  "

$ws.Cells.Item(43, 4).Value = "'3.09"
$ws.Cells.Item(43, 5).Value = "  +13.32%  "

$ws.Cells.Item(44, 4).Value = "'140.61"
$ws.Cells.Item(44, 5).Value = "  -5.38%  "

$ws.Cells.Item(45, 2).Value = "LidoDAOToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(45, 4).Value = "'3.11"
$ws.Cells.Item(45, 5).Value = "  -5.22%  "

$ws.Cells.Item(46, 4).Value = "'2.73"
$ws.Cells.Item(46, 5).Value = "  -10.28%  "

$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value = "'25.08"
$ws.Cells.Item(47, 5).Value = "  +16.12%  "

$ws.Cells.Item(48, 5).Value = "  -1.66%  "

$ws.Cells.Item(49, 2).Value = "WEMIXToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(49, 4).Value = "'2.48"
$ws.Cells.Item(49, 5).Value = "  -9.37%  "

$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(50, 4).Value = "'4.07"
$ws.Cells.Item(50, 5).Value = "  -6.16%  "

$ws.Cells.Item(51, 4).Value = "'0.278"
$ws.Cells.Item(51, 5).Value = "  -10.08%  "

Write-Output "Updated cryptos list"